$wb = $excel.ActiveWorkbook

# --- Resumen sheet ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Cells.Item(2, 2).Value = "Z4"
$wsResumen.Cells.Item(2, 3).Value = 541.0034165855768

# --- Solucion sheet: reshuffle Pedido/Salida pairs ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsSolucion.Cells.Item(2, 1).Value = "Pedido_68"
$wsSolucion.Cells.Item(2, 2).Value = "S001"
$wsSolucion.Cells.Item(3, 1).Value = "Pedido_58"
$wsSolucion.Cells.Item(3, 2).Value = "S021"
$wsSolucion.Cells.Item(4, 1).Value = "Pedido_16"
$wsSolucion.Cells.Item(4, 2).Value = "S041"
$wsSolucion.Cells.Item(5, 1).Value = "Pedido_2"
$wsSolucion.Cells.Item(5, 2).Value = "S061"
$wsSolucion.Cells.Item(6, 1).Value = "Pedido_9"
$wsSolucion.Cells.Item(6, 2).Value = "S051"
$wsSolucion.Cells.Item(7, 1).Value = "Pedido_37"
$wsSolucion.Cells.Item(7, 2).Value = "S011"
$wsSolucion.Cells.Item(8, 1).Value = "Pedido_61"
$wsSolucion.Cells.Item(8, 2).Value = "S031"
$wsSolucion.Cells.Item(9, 1).Value = "Pedido_45"
$wsSolucion.Cells.Item(9, 2).Value = "S022"
$wsSolucion.Cells.Item(10, 1).Value = "Pedido_23"
$wsSolucion.Cells.Item(10, 2).Value = "S042"
$wsSolucion.Cells.Item(11, 1).Value = "Pedido_6"
$wsSolucion.Cells.Item(11, 2).Value = "S071"
$wsSolucion.Cells.Item(12, 1).Value = "Pedido_63"
$wsSolucion.Cells.Item(12, 2).Value = "S002"
$wsSolucion.Cells.Item(13, 1).Value = "Pedido_76"
$wsSolucion.Cells.Item(13, 2).Value = "S012"
$wsSolucion.Cells.Item(14, 1).Value = "Pedido_32"
$wsSolucion.Cells.Item(14, 2).Value = "S062"
$wsSolucion.Cells.Item(15, 1).Value = "Pedido_12"
$wsSolucion.Cells.Item(15, 2).Value = "S052"
$wsSolucion.Cells.Item(16, 1).Value = "Pedido_39"
$wsSolucion.Cells.Item(16, 2).Value = "S032"
$wsSolucion.Cells.Item(17, 1).Value = "Pedido_8"
$wsSolucion.Cells.Item(17, 2).Value = "S003"
$wsSolucion.Cells.Item(18, 1).Value = "Pedido_71"
$wsSolucion.Cells.Item(18, 2).Value = "S072"
$wsSolucion.Cells.Item(19, 1).Value = "Pedido_40"
$wsSolucion.Cells.Item(19, 2).Value = "S023"
$wsSolucion.Cells.Item(20, 1).Value = "Pedido_35"
$wsSolucion.Cells.Item(20, 2).Value = "S043"
$wsSolucion.Cells.Item(21, 1).Value = "Pedido_43"
$wsSolucion.Cells.Item(21, 2).Value = "S013"
$wsSolucion.Cells.Item(22, 1).Value = "Pedido_15"
$wsSolucion.Cells.Item(22, 2).Value = "S063"
$wsSolucion.Cells.Item(23, 1).Value = "Pedido_13"
$wsSolucion.Cells.Item(23, 2).Value = "S033"
$wsSolucion.Cells.Item(24, 1).Value = "Pedido_65"
$wsSolucion.Cells.Item(24, 2).Value = "S053"
$wsSolucion.Cells.Item(25, 1).Value = "Pedido_42"
$wsSolucion.Cells.Item(25, 2).Value = "S004"
$wsSolucion.Cells.Item(26, 1).Value = "Pedido_33"
$wsSolucion.Cells.Item(26, 2).Value = "S044"
$wsSolucion.Cells.Item(27, 1).Value = "Pedido_50"
$wsSolucion.Cells.Item(27, 2).Value = "S073"
$wsSolucion.Cells.Item(28, 1).Value = "Pedido_31"
$wsSolucion.Cells.Item(28, 2).Value = "S024"
$wsSolucion.Cells.Item(29, 1).Value = "Pedido_80"
$wsSolucion.Cells.Item(29, 2).Value = "S014"
$wsSolucion.Cells.Item(30, 1).Value = "Pedido_26"
$wsSolucion.Cells.Item(30, 2).Value = "S064"
$wsSolucion.Cells.Item(31, 1).Value = "Pedido_14"
$wsSolucion.Cells.Item(31, 2).Value = "S005"
$wsSolucion.Cells.Item(32, 1).Value = "Pedido_19"
$wsSolucion.Cells.Item(32, 2).Value = "S034"
$wsSolucion.Cells.Item(33, 1).Value = "Pedido_77"
$wsSolucion.Cells.Item(33, 2).Value = "S054"
$wsSolucion.Cells.Item(34, 1).Value = "Pedido_34"
$wsSolucion.Cells.Item(34, 2).Value = "S074"
$wsSolucion.Cells.Item(35, 1).Value = "Pedido_70"
$wsSolucion.Cells.Item(35, 2).Value = "S045"
$wsSolucion.Cells.Item(36, 1).Value = "Pedido_74"
$wsSolucion.Cells.Item(36, 2).Value = "S015"
$wsSolucion.Cells.Item(37, 1).Value = "Pedido_56"
$wsSolucion.Cells.Item(37, 2).Value = "S025"
$wsSolucion.Cells.Item(38, 1).Value = "Pedido_25"
$wsSolucion.Cells.Item(38, 2).Value = "S055"
$wsSolucion.Cells.Item(39, 1).Value = "Pedido_22"
$wsSolucion.Cells.Item(39, 2).Value = "S065"
$wsSolucion.Cells.Item(40, 1).Value = "Pedido_41"
$wsSolucion.Cells.Item(40, 2).Value = "S035"
$wsSolucion.Cells.Item(41, 1).Value = "Pedido_3"
$wsSolucion.Cells.Item(41, 2).Value = "S006"
$wsSolucion.Cells.Item(42, 1).Value = "Pedido_51"
$wsSolucion.Cells.Item(42, 2).Value = "S046"
$wsSolucion.Cells.Item(43, 1).Value = "Pedido_54"
$wsSolucion.Cells.Item(43, 2).Value = "S026"
$wsSolucion.Cells.Item(44, 1).Value = "Pedido_67"
$wsSolucion.Cells.Item(44, 2).Value = "S075"
$wsSolucion.Cells.Item(45, 1).Value = "Pedido_24"
$wsSolucion.Cells.Item(45, 2).Value = "S056"
$wsSolucion.Cells.Item(46, 1).Value = "Pedido_47"
$wsSolucion.Cells.Item(46, 2).Value = "S036"
$wsSolucion.Cells.Item(47, 1).Value = "Pedido_49"
$wsSolucion.Cells.Item(47, 2).Value = "S066"
$wsSolucion.Cells.Item(48, 1).Value = "Pedido_21"
$wsSolucion.Cells.Item(48, 2).Value = "S016"
$wsSolucion.Cells.Item(49, 1).Value = "Pedido_20"
$wsSolucion.Cells.Item(49, 2).Value = "S047"
$wsSolucion.Cells.Item(50, 1).Value = "Pedido_30"
$wsSolucion.Cells.Item(50, 2).Value = "S027"
$wsSolucion.Cells.Item(51, 1).Value = "Pedido_62"
$wsSolucion.Cells.Item(51, 2).Value = "S076"
$wsSolucion.Cells.Item(52, 1).Value = "Pedido_44"
$wsSolucion.Cells.Item(52, 2).Value = "S007"
$wsSolucion.Cells.Item(53, 1).Value = "Pedido_36"
$wsSolucion.Cells.Item(53, 2).Value = "S037"
$wsSolucion.Cells.Item(54, 1).Value = "Pedido_69"
$wsSolucion.Cells.Item(54, 2).Value = "S067"
$wsSolucion.Cells.Item(55, 1).Value = "Pedido_60"
$wsSolucion.Cells.Item(55, 2).Value = "S057"
$wsSolucion.Cells.Item(56, 1).Value = "Pedido_1"
$wsSolucion.Cells.Item(56, 2).Value = "S017"
$wsSolucion.Cells.Item(57, 1).Value = "Pedido_18"
$wsSolucion.Cells.Item(57, 2).Value = "S077"
$wsSolucion.Cells.Item(58, 1).Value = "Pedido_75"
$wsSolucion.Cells.Item(58, 2).Value = "S048"
$wsSolucion.Cells.Item(59, 1).Value = "Pedido_52"
$wsSolucion.Cells.Item(59, 2).Value = "S028"
$wsSolucion.Cells.Item(60, 1).Value = "Pedido_28"
$wsSolucion.Cells.Item(60, 2).Value = "S008"
$wsSolucion.Cells.Item(61, 1).Value = "Pedido_78"
$wsSolucion.Cells.Item(61, 2).Value = "S068"
$wsSolucion.Cells.Item(62, 1).Value = "Pedido_55"
$wsSolucion.Cells.Item(62, 2).Value = "S058"
$wsSolucion.Cells.Item(63, 1).Value = "Pedido_7"
$wsSolucion.Cells.Item(63, 2).Value = "S038"
$wsSolucion.Cells.Item(64, 1).Value = "Pedido_5"
$wsSolucion.Cells.Item(64, 2).Value = "S018"
$wsSolucion.Cells.Item(65, 1).Value = "Pedido_4"
$wsSolucion.Cells.Item(65, 2).Value = "S078"
$wsSolucion.Cells.Item(66, 1).Value = "Pedido_11"
$wsSolucion.Cells.Item(66, 2).Value = "S049"
$wsSolucion.Cells.Item(67, 1).Value = "Pedido_29"
$wsSolucion.Cells.Item(67, 2).Value = "S029"
$wsSolucion.Cells.Item(68, 1).Value = "Pedido_27"
$wsSolucion.Cells.Item(68, 2).Value = "S069"
$wsSolucion.Cells.Item(69, 1).Value = "Pedido_48"
$wsSolucion.Cells.Item(69, 2).Value = "S059"
$wsSolucion.Cells.Item(70, 1).Value = "Pedido_59"
$wsSolucion.Cells.Item(70, 2).Value = "S009"
$wsSolucion.Cells.Item(71, 1).Value = "Pedido_73"
$wsSolucion.Cells.Item(71, 2).Value = "S079"
$wsSolucion.Cells.Item(72, 1).Value = "Pedido_66"
$wsSolucion.Cells.Item(72, 2).Value = "S039"
$wsSolucion.Cells.Item(73, 1).Value = "Pedido_53"
$wsSolucion.Cells.Item(73, 2).Value = "S050"
$wsSolucion.Cells.Item(74, 1).Value = "Pedido_38"
$wsSolucion.Cells.Item(74, 2).Value = "S030"
$wsSolucion.Cells.Item(75, 1).Value = "Pedido_46"
$wsSolucion.Cells.Item(75, 2).Value = "S070"
$wsSolucion.Cells.Item(76, 1).Value = "Pedido_17"
$wsSolucion.Cells.Item(76, 2).Value = "S019"
$wsSolucion.Cells.Item(77, 1).Value = "Pedido_10"
$wsSolucion.Cells.Item(77, 2).Value = "S060"
$wsSolucion.Cells.Item(78, 1).Value = "Pedido_64"
$wsSolucion.Cells.Item(78, 2).Value = "S040"
$wsSolucion.Cells.Item(79, 1).Value = "Pedido_79"
$wsSolucion.Cells.Item(79, 2).Value = "S080"
$wsSolucion.Cells.Item(80, 1).Value = "Pedido_57"
$wsSolucion.Cells.Item(80, 2).Value = "S010"
$wsSolucion.Cells.Item(81, 1).Value = "Pedido_72"
$wsSolucion.Cells.Item(81, 2).Value = "S020"

# --- Metricas sheet ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Cells.Item(2, 2).Value = 539.4052113742026
$wsMetricas.Cells.Item(3, 2).Value = 539.7534328035464
$wsMetricas.Cells.Item(4, 2).Value = 536.1188885284896
$wsMetricas.Cells.Item(5, 2).Value = 541.0034165855768
